# Update the "Price" column (D) values for the symbol list refresh,
# per the "Updated symbol list on Wed Dec 14 05:55:17 UTC 2022 with GitHub Actions" commit.
#
# These cells are stored as text (inline strings) in the original workbook, so we
# prefix the assigned values with a leading apostrophe to force Excel to keep them
# as text instead of auto-coercing the numeric-looking strings into numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "'274.81"
$ws.Range("D3").Value  = "'22.91"
$ws.Range("D4").Value  = "'6.403"
$ws.Range("D6").Value  = "'3.670"
$ws.Range("D7").Value  = "'6.667"
$ws.Range("D8").Value  = "'1.398"
$ws.Range("D9").Value  = "'0.8328"
$ws.Range("D10").Value = "'0.01376"
$ws.Range("D12").Value = "'0.08307"
$ws.Range("D14").Value = "'0.03099"
$ws.Range("D15").Value = "'0.09300"
$ws.Range("D16").Value = "'3.844"
$ws.Range("D17").Value = "'0.001636"
$ws.Range("D18").Value = "'0.04781"
$ws.Range("D19").Value = "'0.006371"
$ws.Range("D20").Value = "'0.005671"
$ws.Range("D22").Value = "'0.0001501"
$ws.Range("D23").Value = "'3.714"
$ws.Range("D27").Value = "'0.0002681"
$ws.Range("D40").Value = "'0.04709"
$ws.Range("D41").Value = "'0.007049"
$ws.Range("D42").Value = "'0.1159"
$ws.Range("D43").Value = "'0.003702"
$ws.Range("D44").Value = "'0.01181"
$ws.Range("D45").Value = "'0.00006264"
$ws.Range("D48").Value = "'0.7969"
$ws.Range("D49").Value = "'0.03820"
$ws.Range("D50").Value = "'0.00002301"
$ws.Range("D51").Value = "'0.01241"
